$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 (RM 232) entirely - remaining rows shift up by one
$ws.Rows("26:26").Delete()

# After that delete, "SC 92" (originally row 28) is now row 27 - delete it too
$ws.Rows("27:27").Delete()

# Fix up the C column values that moved/changed as part of this edit
# Row 26 is now "SC 5": C becomes 10.8
$ws.Range("C26").Value = 10.8

# Row 27 is now "SC 101": C becomes blank
$ws.Range("C27").Value = ""

# Row 33 is now "SC 232": C becomes 10.4
$ws.Range("C33").Value = 10.4
